# Update "想去人数" (column F) values for several rows across the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets, plus one row on
# "演出" (sheet2), matching the refreshed data snapshot from the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F16").Value = 6660
$ws1.Range("F20").Value = 2172
$ws1.Range("F21").Value = 2997
$ws1.Range("F23").Value = 202
$ws1.Range("F25").Value = 1701
$ws1.Range("F27").Value = 282
$ws1.Range("F37").Value = 883
$ws1.Range("F39").Value = 418

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 23

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F20").Value = 23
$ws4.Range("F23").Value = 6661
$ws4.Range("F26").Value = 2172
$ws4.Range("F27").Value = 2997
$ws4.Range("F30").Value = 202
$ws4.Range("F33").Value = 1701
$ws4.Range("F36").Value = 282
$ws4.Range("F46").Value = 883
$ws4.Range("F48").Value = 418
